$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.05247427851054
$ws.Range("C2").Value = 5.481286445270841
$ws.Range("D2").Value = 5.395086469004689
$ws.Range("F2").Value = 25.75112018213443
$ws.Range("G2").Value = 31.27472198665229
$ws.Range("H2").Value = 14.989452409211
$ws.Range("K2").Value = 9.285031372711872
$ws.Range("B3").Value = 9.736296559564558
$ws.Range("C3").Value = 5.298933014485848
$ws.Range("D3").Value = 5.34801924787769
$ws.Range("F3").Value = 25.84208167849445
$ws.Range("G3").Value = 31.43040331321058
$ws.Range("H3").Value = 15.05818048246753
$ws.Range("K3").Value = 9.060503429531419
$ws.Range("B4").Value = 9.538590381248603
$ws.Range("C4").Value = 5.182829625030983
$ws.Range("D4").Value = 5.318665119364594
$ws.Range("F4").Value = 25.90627772511543
$ws.Range("G4").Value = 31.53861700147895
$ws.Range("H4").Value = 15.10337583662439
$ws.Range("K4").Value = 8.921228494228563
$ws.Range("B5").Value = 9.45725413334854
$ws.Range("C5").Value = 5.134521825509039
$ws.Range("D5").Value = 5.306594725813303
$ws.Range("F5").Value = 25.93452639921431
$ws.Range("G5").Value = 31.58586655027319
$ws.Range("H5").Value = 15.12254608537365
$ws.Range("K5").Value = 8.864200138320705
$ws.Range("B6").Value = 9.44370579311625
$ws.Range("C6").Value = 5.126441703009845
$ws.Range("D6").Value = 5.304584077921636
$ws.Range("F6").Value = 25.9393428865236
$ws.Range("G6").Value = 31.59390196500571
$ws.Range("H6").Value = 15.12577473919929
$ws.Range("K6").Value = 8.854716644134282
$ws.Range("B7").Value = 9.537496390158577
$ws.Range("C7").Value = 5.182182092845555
$ws.Range("D7").Value = 5.318502764588338
$ws.Range("F7").Value = 25.90665025525917
$ws.Range("G7").Value = 31.53924149567582
$ws.Range("H7").Value = 15.10363132614031
$ws.Range("K7").Value = 8.920460381901691
$ws.Range("B8").Value = 9.944271803123703
$ws.Range("C8").Value = 5.41929571697543
$ws.Range("D8").Value = 5.378955625253218
$ws.Range("F8").Value = 25.78074564092947
$ws.Range("G8").Value = 31.32576702602648
$ws.Range("H8").Value = 15.01252765470338
$ws.Range("K8").Value = 9.20795614205486
$ws.Range("B9").Value = 10.70845377239852
$ws.Range("C9").Value = 5.849672665236339
$ws.Range("D9").Value = 5.493637742585423
$ws.Range("F9").Value = 25.60049667375626
$ws.Range("G9").Value = 31.0083097569761
$ws.Range("H9").Value = 14.85767302266133
$ws.Range("K9").Value = 9.757073078695228
$ws.Range("B10").Value = 11.24326151388055
$ws.Range("C10").Value = 6.142713971452904
$ws.Range("D10").Value = 5.575199484428476
$ws.Range("F10").Value = 25.50925383793213
$ws.Range("G10").Value = 30.83805143289617
$ws.Range("H10").Value = 14.75844761351873
$ws.Range("K10").Value = 10.14726716744912
$ws.Range("B11").Value = 11.47965940838266
$ws.Range("C11").Value = 6.270642730555852
$ws.Range("D11").Value = 5.611648809808665
$ws.Range("F11").Value = 25.47679345659133
$ws.Range("G11").Value = 30.77453140056531
$ws.Range("H11").Value = 14.71647469132315
$ws.Range("K11").Value = 10.32108469942936
$ws.Range("B12").Value = 11.56810622769365
$ws.Range("C12").Value = 6.318287888524289
$ws.Range("D12").Value = 5.625351464098494
$ws.Range("F12").Value = 25.4658099045444
$ws.Range("G12").Value = 30.75250106917444
$ws.Range("H12").Value = 14.7010366459121
$ws.Range("K12").Value = 10.38631376977099
$ws.Range("B13").Value = 11.54910649815009
$ws.Range("C13").Value = 6.308062543976082
$ws.Range("D13").Value = 5.622404889160364
$ws.Range("F13").Value = 25.46811710058552
$ws.Range("G13").Value = 30.75715540138747
$ws.Range("H13").Value = 14.70434119925043
$ws.Range("K13").Value = 10.37229278196218
$ws.Range("B14").Value = 11.4869578989368
$ws.Range("C14").Value = 6.274578667510731
$ws.Range("D14").Value = 5.612778171392988
$ws.Range("F14").Value = 25.47586356821092
$ws.Range("G14").Value = 30.7726782911127
$ws.Range("H14").Value = 14.71519544528235
$ws.Range("K14").Value = 10.32646333425986
$ws.Range("B15").Value = 11.44874822789951
$ws.Range("C15").Value = 6.253964113915552
$ws.Range("D15").Value = 5.606868343046157
$ws.Range("F15").Value = 25.48077912010777
$ws.Range("G15").Value = 30.78245057961457
$ws.Range("H15").Value = 14.72190342495775
$ws.Range("K15").Value = 10.29831263766028
$ws.Range("B16").Value = 11.22766657842133
$ws.Range("C16").Value = 6.134243072363174
$ws.Range("D16").Value = 5.572803774121843
$ws.Range("F16").Value = 25.51155790717885
$ws.Range("G16").Value = 30.84248468880654
$ws.Range("H16").Value = 14.76125442268695
$ws.Range("K16").Value = 10.13582823128608
$ws.Range("B17").Value = 11.09021503222192
$ws.Range("C17").Value = 6.059401697787274
$ws.Range("D17").Value = 5.551734615426523
$ws.Range("F17").Value = 25.53276237726179
$ws.Range("G17").Value = 30.8828964365842
$ws.Range("H17").Value = 14.78620639926722
$ws.Range("K17").Value = 10.03516002549207
$ws.Range("B18").Value = 11.01051240428533
$ws.Range("C18").Value = 6.01585009269162
$ws.Range("D18").Value = 5.53955502743107
$ws.Range("F18").Value = 25.54580971391819
$ws.Range("G18").Value = 30.9074503630844
$ws.Range("H18").Value = 14.80085597722999
$ws.Range("K18").Value = 9.976915005377133
$ws.Range("B19").Value = 10.98341845831779
$ws.Range("C19").Value = 6.001018392841844
$ws.Range("D19").Value = 5.535420906143164
$ws.Range("F19").Value = 25.5503732478433
$ws.Range("G19").Value = 30.91598829232175
$ws.Range("H19").Value = 14.80586719951334
$ws.Range("K19").Value = 9.957137250129769
$ws.Range("B20").Value = 11.10491426681603
$ws.Range("C20").Value = 6.067421132216066
$ws.Range("D20").Value = 5.553983837133159
$ws.Range("F20").Value = 25.53041698325598
$ws.Range("G20").Value = 30.87845879176639
$ws.Range("H20").Value = 14.7835193838508
$ws.Range("K20").Value = 10.04591232723894
$ws.Range("B21").Value = 11.50524214522723
$ws.Range("C21").Value = 6.284435556583788
$ws.Range("D21").Value = 5.615608530060275
$ws.Range("F21").Value = 25.47355267401286
$ws.Range("G21").Value = 30.76806377608763
$ws.Range("H21").Value = 14.71199490085885
$ws.Range("K21").Value = 10.33994108737776
$ws.Range("B22").Value = 11.76059991824221
$ws.Range("C22").Value = 6.421602091944625
$ws.Range("D22").Value = 5.655298795324633
$ws.Range("F22").Value = 25.44401841308047
$ws.Range("G22").Value = 30.70771615919485
$ws.Range("H22").Value = 14.66790903934999
$ws.Range("K22").Value = 10.52863075342591
$ws.Range("B23").Value = 11.62490986198037
$ws.Range("C23").Value = 6.348828156435057
$ws.Range("D23").Value = 5.634170811182677
$ws.Range("F23").Value = 25.4590809591585
$ws.Range("G23").Value = 30.7388386462416
$ws.Range("H23").Value = 14.69119479668949
$ws.Range("K23").Value = 10.42826075024451
$ws.Range("B24").Value = 11.09827085433188
$ws.Range("C24").Value = 6.063797177065802
$ws.Range("D24").Value = 5.552967171131562
$ws.Range("F24").Value = 25.53147466817851
$ws.Range("G24").Value = 30.88046093866169
$ws.Range("H24").Value = 14.7847332355265
$ws.Range("K24").Value = 10.04105235783738
$ws.Range("B25").Value = 10.505989124921
$ws.Range("C25").Value = 5.737177739532009
$ws.Range("D25").Value = 5.46306232877331
$ws.Range("F25").Value = 25.64206368456227
$ws.Range("G25").Value = 31.08322172128228
$ws.Range("H25").Value = 14.89701365656715
$ws.Range("K25").Value = 9.610554420126002
